$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LogIn")
$ws.Range("B1").Value = "admin"
